$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")

# Row 5 - Compilation success
$ws.Range("B5").Value = "no"
$ws.Range("C5").Value = "Called wrong method"

# Row 6 - Runtime without error
$ws.Range("B6").ClearContents()

# Row 7 - Assertion validity
$ws.Range("B7").ClearContents()
$ws.Range("C7").ClearContents()

# Row 12 - Code BLEU
$ws.Range("B12").Value = 0.3106136319912257
$ws.Range("C12").Value = "{'codebleu': 0.31061363199122566, 'ngram_match_score': 0.0704367794283077, 'weighted_ngram_match_score': 0.09453205847935507, 'syntax_match_score': 0.5943396226415094, 'dataflow_match_score': 0.48314606741573035}"

# Update selection to C6 on the active sheet
$ws.Activate()
$ws.Range("C6").Select()
